$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Brandyn Offley" -> split into "Brandyn" / " " / "Offley" runs,
#        each name wrapped in a spellStart/spellEnd proofErr pair.
$p1 = $d.Paragraphs(1).Range
$xmlBrandyn = '<w:p ' + $W + '>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Brandyn</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Offley</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'
$p1.InsertXML($xmlBrandyn) | Out-Null

# --- 2. "Ira Grunwell" -> split into "Ira " / "Grunwell" runs, with only
#        "Grunwell" wrapped in a spellStart/spellEnd proofErr pair.
$p5 = $d.Paragraphs(5).Range
$xmlIra = '<w:p ' + $W + '>' + `
    '<w:r><w:t xml:space="preserve">Ira </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Grunwell</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'
$p5.InsertXML($xmlIra) | Out-Null

# --- 3. Replace the last paragraph ("Cell: 778-866-6987", which also
#        carries the trailing _GoBack bookmark) with a clean copy of
#        itself (bookmark dropped here, re-added at the very end below)
#        followed by the new Ryan Dieno contact block:
#          <blank paragraph>
#          Ryan Dieno                     (Dieno wrapped in proofErr)
#          <blank placeholder paragraph>  (filled in with the hyperlink below)
#          Cell: 778-322-2349
$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIdx).Range
$xmlBlock = '<w:p ' + $W + '><w:r><w:t>Cell: 778-866-6987</w:t></w:r></w:p>' + `
    '<w:p ' + $W + '/>' + `
    '<w:p ' + $W + '>' + `
        '<w:r><w:t xml:space="preserve">Ryan </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>Dieno</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>' + `
    '<w:p ' + $W + '/>' + `
    '<w:p ' + $W + '><w:r><w:t>Cell: 778-322-2349</w:t></w:r></w:p>'
$pLast.InsertXML($xmlBlock) | Out-Null

# --- 4. Fill the placeholder paragraph with "Email: " + the new hyperlink
$emailIdx = $lastIdx + 3
$pEmail = $d.Paragraphs($emailIdx).Range
$d.Hyperlinks.Add($pEmail, "mailto:ryan.dieno@gmail.com", $null, $null, "ryan.dieno@gmail.com") | Out-Null
$pEmailStart = $d.Paragraphs($emailIdx).Range
$insPt = $d.Range($pEmailStart.Start, $pEmailStart.Start)
$insPt.InsertBefore("Email: ") | Out-Null

# --- 5. Re-add the _GoBack bookmark on the new final paragraph
$finalIdx = $d.Paragraphs.Count
$pFinal = $d.Paragraphs($finalIdx).Range
$xmlFinal = '<w:p ' + $W + '><w:r><w:t>Cell: 778-322-2349</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'
$pFinal.InsertXML($xmlFinal) | Out-Null

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
